# Update "Glosario Data Climático.xlsx"
# The "Definición" column (C) entries each start with a leading space followed
# by a lower-case word (e.g. " gas de efecto invernadero...").
# This edit capitalizes the first letter of the actual text in each of those
# definition cells (rows 4 through 21), leaving the leading space intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 4; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $text = $cell.Value()

    if ($text -ne $null -and $text -ne "") {
        # Find index of first non-space character
        $trimmed = $text.TrimStart(" ")
        $leadLen = $text.Length - $trimmed.Length
        $lead = $text.Substring(0, $leadLen)

        if ($trimmed.Length -gt 0) {
            $newFirst = $trimmed.Substring(0, 1).ToUpper()
            $rest = $trimmed.Substring(1)
            $cell.Value = $lead + $newFirst + $rest
        }
    }
}

# Restore the selected cell as recorded in the saved workbook view.
$ws.Range("C22").Select()
